$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.653.47"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "2.640.57"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "2.673.80"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "3.130.30"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "59.472.70"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "2.662.39"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.429"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "2.755.56"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "0.0₃0845"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.64%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +18.50%  "
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.894"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "293.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.630"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.992"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0550"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0235"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
